# Update gh-pages output with latest scraped event stats (commit 456a3b4).
# Two worksheets ("展览" and "全部类型") carry duplicate rows for the same
# events, so the "想去人数" (F column) figures need bumping in both.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — rows 2,4,7,8
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5562
$ws1.Range("F4").Value = 638
$ws1.Range("F7").Value = 46
$ws1.Range("F8").Value = 366

# Sheet "全部类型" (all types) — same events, but shifted down one row
# (row 8 there holds an extra event absent from "展览")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5562
$ws4.Range("F4").Value = 638
$ws4.Range("F7").Value = 46
$ws4.Range("F9").Value = 366
